# Append a new data row (row 89) to each of the four worksheets, mirroring
# the most recent existing row (row 88) but with the next day's timestamp
# and refreshed payload/checksum values (security-vulnerability-check log
# update).
$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = 1
        A = 45875.43287037037
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x30"
        E = "0x14"
        F = 380
        G = 759863127514710945038336.0
        H = 304
        I = 14
    },
    @{
        Sheet = 2
        A = 45875.43287037037
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x34"
        E = "0xe"
        F = 380
        G = 568432987514711010443264.0
        H = 308
        I = 14
    },
    @{
        Sheet = 3
        A = 45875.43287037037
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x76"
        E = "0x7"
        F = 130
        G = 568631262647113970876416.0
        H = 118
        I = 7
    },
    @{
        Sheet = 4
        A = 45875.43287037037
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x75"
        E = "0x3"
        F = 130
        G = 985046333984776009023488.0
        H = 117
        I = 3
    }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)
    $newRow = 89

    $ws.Cells.Item($newRow, 1).Value = $r.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $r.B
    $ws.Cells.Item($newRow, 3).Value = $r.C
    $ws.Cells.Item($newRow, 4).Value = $r.D
    $ws.Cells.Item($newRow, 5).Value = $r.E

    $ws.Cells.Item($newRow, 6).Value = $r.F
    $ws.Cells.Item($newRow, 7).Value = $r.G
    $ws.Cells.Item($newRow, 8).Value = $r.H
    $ws.Cells.Item($newRow, 9).Value = $r.I
}

Write-Host "Appended row 89 to" $rows.Count "sheets"
